# Applies a data-shuffle edit to the "Scotland Premiership" sheet.
# For a number of small groups of rows that all share the same Date (column D),
# the match-specific data (everything except the running id in column A,
# the Div in column C and the Date in column D) is redistributed among the
# rows of the group according to a fixed permutation. This script snapshots
# the "moving" data (column B, and the block E:AD) for every row involved,
# then writes the snapshots back to the rows in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together with a match (everything but id/Div/Date):
#   B            = id number (match id, NOT the running row id in col A)
#   E .. AD      = HomeTeam .. PL_AhUnder  (columns 5 .. 30)
$dataColumns = @(2) + @(5..30)

function Get-RowSnapshot($ws, $row) {
    $vals = @{}
    foreach ($c in $dataColumns) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowSnapshot($ws, $row, $snapshot) {
    foreach ($c in $dataColumns) {
        $ws.Cells.Item($row, $c).Value = $snapshot[$c]
    }
}

function Apply-Permutation($ws, [int[]]$rows, [int[]]$sourceRows) {
    # Take a snapshot of every row first (before any writes happen),
    # then write snapshot[sourceRows[i]] into rows[i].
    $snapshots = @{}
    foreach ($r in $rows) {
        $snapshots[$r] = Get-RowSnapshot $ws $r
    }
    for ($i = 0; $i -lt $rows.Count; $i++) {
        Set-RowSnapshot $ws $rows[$i] $snapshots[$sourceRows[$i]]
    }
}

# Group: rows 61-63 (rotate up by one, wrapping)
Apply-Permutation $ws @(61,62,63) @(62,63,61)

# Group: rows 76-77 (swap)
Apply-Permutation $ws @(76,77) @(77,76)

# Group: rows 79-80 (swap)
Apply-Permutation $ws @(79,80) @(80,79)

# Group: rows 128-132 (rotate up by one, wrapping)
Apply-Permutation $ws @(128,129,130,131,132) @(129,130,131,132,128)

# Group: rows 133-134 (swap)
Apply-Permutation $ws @(133,134) @(134,133)

# Group: rows 163-166 (full reversal)
Apply-Permutation $ws @(163,164,165,166) @(166,165,164,163)

# Group: rows 176-177 (swap)
Apply-Permutation $ws @(176,177) @(177,176)

# Group: rows 179-180 (swap)
Apply-Permutation $ws @(179,180) @(180,179)

# Group: rows 182 & 184 (swap; row 183 is untouched)
Apply-Permutation $ws @(182,184) @(184,182)

# Group: rows 227-228 (swap)
Apply-Permutation $ws @(227,228) @(228,227)
